{"js": "// Append the new \"DATA STORAGE\" section to the end of the document body,\n// right after the paragraph that ends the SECURITY PLAN section\n// (\"...the device works offline and syncs later.)\").\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst anchorText = \"(If there\\u2019s no internet, the device works offline and syncs later.)\";\nlet anchorPara = null;\nfor (let i = paragraphs.items.length - 1; i >= 0; i--) {\n  if (paragraphs.items[i].text.indexOf(anchorText) !== -1) {\n    anchorPara = paragraphs.items[i];\n    break;\n  }\n}\n// Fall back to the very last paragraph of the body if the text wasn't found\n// (keeps the script resilient to minor whitespace/quote differences).\nif (!anchorPara) {\n  anchorPara = paragraphs.items[paragraphs.items.length - 1];\n}\n\n// The new paragraphs to append, in order. Empty strings become blank lines.\nconst newParagraphTexts = [\n  \"\",\n  \"DATA STORAGE \",\n  \"What data the device will gather\",\n  \"\",\n  \"From the Pi (device-side):\",\n  \"\",\n  \"Reminder events: when a scheduled alert fires (timestamp).\",\n  \"Sensor events: when the motor/buzzer/LED is activated and for how long.\",\n  \"User actions: \\u201cTaken\\u201d(Open/closed Lid), \\u201cSnooze\\u201d, \\u201cDismiss\\u201d (from the phone app command).\",\n  \"Battery telemetry: percentage/voltage (via PiSugar API), charging state, last-low-battery time. (Link to documentation for this: https://www.pisugar.com/blogs/pisugar-blog/display-raspberry-pi-battery-indicator)\",\n  \"\",\n  \"Device health: uptime, last sync time.\",\n];\n\nlet insertAfter = anchorPara;\nfor (const text of newParagraphTexts) {\n  insertAfter = insertAfter.insertParagraph(text, Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "# Append the new \"DATA STORAGE\" section to the end of the document,\n# right after the paragraph that ends the SECURITY PLAN section\n# (\"...the device works offline and syncs later.)\").\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that ends with the SECURITY PLAN section's closing\n# remark, searching from the end of the document for resilience.\n$count = $d.Paragraphs.Count\n$anchorIndex = $count\nfor ($i = $count; $i -ge 1; $i--) {\n  if ($d.Paragraphs.Item($i).Range.Text -like \"*syncs later.)*\") {\n    $anchorIndex = $i\n    break\n  }\n}\n$anchor = $d.Paragraphs.Item($anchorIndex)\n\n# New paragraphs to append, in order. Empty strings become blank lines.\n$newParagraphTexts = @(\n  \"\",\n  \"DATA STORAGE \",\n  \"What data the device will gather\",\n  \"\",\n  \"From the Pi (device-side):\",\n  \"\",\n  \"Reminder events: when a scheduled alert fires (timestamp).\",\n  \"Sensor events: when the motor/buzzer/LED is activated and for how long.\",\n  \"User actions: \u201cTaken\u201d(Open/closed Lid), \u201cSnooze\u201d, \u201cDismiss\u201d (from the phone app command).\",\n  \"Battery telemetry: percentage/voltage (via PiSugar API), charging state, last-low-battery time. (Link to documentation for this: https://www.pisugar.com/blogs/pisugar-blog/display-raspberry-pi-battery-indicator)\",\n  \"\",\n  \"Device health: uptime, last sync time.\"\n)\n\nforeach ($text in $newParagraphTexts) {\n  $r = $anchor.Range\n  $r.Collapse(0)\n  $r.InsertParagraphAfter()\n  $anchor = $d.Paragraphs.Item($d.Paragraphs.Count)\n  if ($text -ne \"\") {\n    $anchor.Range.Text = $text\n  }\n}\n"}
